$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Timestamp (A), Actual Consumption MW (B), and Lookup (D) columns
# for rows 2-193. Column C (Quarter) is unchanged.

$ws.Cells.Item(2, 1).Value = 45796
$ws.Cells.Item(2, 2).Value = 4833
$ws.Cells.Item(2, 4).Value = "19.05.20251"
$ws.Cells.Item(3, 1).Value = 45796.01041666666
$ws.Cells.Item(3, 2).Value = 4804
$ws.Cells.Item(3, 4).Value = "19.05.20252"
$ws.Cells.Item(4, 1).Value = 45796.02083333334
$ws.Cells.Item(4, 2).Value = 4811
$ws.Cells.Item(4, 4).Value = "19.05.20253"
$ws.Cells.Item(5, 1).Value = 45796.03125
$ws.Cells.Item(5, 2).Value = 4734
$ws.Cells.Item(5, 4).Value = "19.05.20254"
$ws.Cells.Item(6, 1).Value = 45796.04166666666
$ws.Cells.Item(6, 2).Value = 4730
$ws.Cells.Item(6, 4).Value = "19.05.20255"
$ws.Cells.Item(7, 1).Value = 45796.05208333334
$ws.Cells.Item(7, 2).Value = 4675
$ws.Cells.Item(7, 4).Value = "19.05.20256"
$ws.Cells.Item(8, 1).Value = 45796.0625
$ws.Cells.Item(8, 2).Value = 4653
$ws.Cells.Item(8, 4).Value = "19.05.20257"
$ws.Cells.Item(9, 1).Value = 45796.07291666666
$ws.Cells.Item(9, 2).Value = 4666
$ws.Cells.Item(9, 4).Value = "19.05.20258"
$ws.Cells.Item(10, 1).Value = 45796.08333333334
$ws.Cells.Item(10, 2).Value = 4653
$ws.Cells.Item(10, 4).Value = "19.05.20259"
$ws.Cells.Item(11, 1).Value = 45796.09375
$ws.Cells.Item(11, 2).Value = 4663
$ws.Cells.Item(11, 4).Value = "19.05.202510"
$ws.Cells.Item(12, 1).Value = 45796.10416666666
$ws.Cells.Item(12, 2).Value = 4641
$ws.Cells.Item(12, 4).Value = "19.05.202511"
$ws.Cells.Item(13, 1).Value = 45796.11458333334
$ws.Cells.Item(13, 2).Value = 4608
$ws.Cells.Item(13, 4).Value = "19.05.202512"
$ws.Cells.Item(14, 1).Value = 45796.125
$ws.Cells.Item(14, 2).Value = 4657
$ws.Cells.Item(14, 4).Value = "19.05.202513"
$ws.Cells.Item(15, 1).Value = 45796.13541666666
$ws.Cells.Item(15, 2).Value = 4656
$ws.Cells.Item(15, 4).Value = "19.05.202514"
$ws.Cells.Item(16, 1).Value = 45796.14583333334
$ws.Cells.Item(16, 2).Value = 4707
$ws.Cells.Item(16, 4).Value = "19.05.202515"
$ws.Cells.Item(17, 1).Value = 45796.15625
$ws.Cells.Item(17, 2).Value = 4749
$ws.Cells.Item(17, 4).Value = "19.05.202516"
$ws.Cells.Item(18, 1).Value = 45796.16666666666
$ws.Cells.Item(18, 2).Value = 4768
$ws.Cells.Item(18, 4).Value = "19.05.202517"
$ws.Cells.Item(19, 1).Value = 45796.17708333334
$ws.Cells.Item(19, 2).Value = 4824
$ws.Cells.Item(19, 4).Value = "19.05.202518"
$ws.Cells.Item(20, 1).Value = 45796.1875
$ws.Cells.Item(20, 2).Value = 4806
$ws.Cells.Item(20, 4).Value = "19.05.202519"
$ws.Cells.Item(21, 1).Value = 45796.19791666666
$ws.Cells.Item(21, 2).Value = 4893
$ws.Cells.Item(21, 4).Value = "19.05.202520"
$ws.Cells.Item(22, 1).Value = 45796.20833333334
$ws.Cells.Item(22, 2).Value = 5124
$ws.Cells.Item(22, 4).Value = "19.05.202521"
$ws.Cells.Item(23, 1).Value = 45796.21875
$ws.Cells.Item(23, 2).Value = 5266
$ws.Cells.Item(23, 4).Value = "19.05.202522"
$ws.Cells.Item(24, 1).Value = 45796.22916666666
$ws.Cells.Item(24, 2).Value = 5328
$ws.Cells.Item(24, 4).Value = "19.05.202523"
$ws.Cells.Item(25, 1).Value = 45796.23958333334
$ws.Cells.Item(25, 2).Value = 5532
$ws.Cells.Item(25, 4).Value = "19.05.202524"
$ws.Cells.Item(26, 1).Value = 45796.25
$ws.Cells.Item(26, 2).Value = 5748
$ws.Cells.Item(26, 4).Value = "19.05.202525"
$ws.Cells.Item(27, 1).Value = 45796.26041666666
$ws.Cells.Item(27, 2).Value = 5745
$ws.Cells.Item(27, 4).Value = "19.05.202526"
$ws.Cells.Item(28, 1).Value = 45796.27083333334
$ws.Cells.Item(28, 2).Value = 5858
$ws.Cells.Item(28, 4).Value = "19.05.202527"
$ws.Cells.Item(29, 1).Value = 45796.28125
$ws.Cells.Item(29, 2).Value = 5920
$ws.Cells.Item(29, 4).Value = "19.05.202528"
$ws.Cells.Item(30, 1).Value = 45796.29166666666
$ws.Cells.Item(30, 2).Value = 6207
$ws.Cells.Item(30, 4).Value = "19.05.202529"
$ws.Cells.Item(31, 1).Value = 45796.30208333334
$ws.Cells.Item(31, 2).Value = 5995
$ws.Cells.Item(31, 4).Value = "19.05.202530"
$ws.Cells.Item(32, 1).Value = 45796.3125
$ws.Cells.Item(32, 2).Value = 5975
$ws.Cells.Item(32, 4).Value = "19.05.202531"
$ws.Cells.Item(33, 1).Value = 45796.32291666666
$ws.Cells.Item(33, 2).Value = 5898
$ws.Cells.Item(33, 4).Value = "19.05.202532"
$ws.Cells.Item(34, 1).Value = 45796.33333333334
$ws.Cells.Item(34, 2).Value = 5781
$ws.Cells.Item(34, 4).Value = "19.05.202533"
$ws.Cells.Item(35, 1).Value = 45796.34375
$ws.Cells.Item(35, 2).Value = 5869
$ws.Cells.Item(35, 4).Value = "19.05.202534"
$ws.Cells.Item(36, 1).Value = 45796.35416666666
$ws.Cells.Item(36, 2).Value = 5855
$ws.Cells.Item(36, 4).Value = "19.05.202535"
$ws.Cells.Item(37, 1).Value = 45796.36458333334
$ws.Cells.Item(37, 2).Value = 5857
$ws.Cells.Item(37, 4).Value = "19.05.202536"
$ws.Cells.Item(38, 1).Value = 45796.375
$ws.Cells.Item(38, 2).Value = 5712
$ws.Cells.Item(38, 4).Value = "19.05.202537"
$ws.Cells.Item(39, 1).Value = 45796.38541666666
$ws.Cells.Item(39, 2).Value = 5596
$ws.Cells.Item(39, 4).Value = "19.05.202538"
$ws.Cells.Item(40, 1).Value = 45796.39583333334
$ws.Cells.Item(40, 2).Value = 5630
$ws.Cells.Item(40, 4).Value = "19.05.202539"
$ws.Cells.Item(41, 1).Value = 45796.40625
$ws.Cells.Item(41, 2).Value = 5583
$ws.Cells.Item(41, 4).Value = "19.05.202540"
$ws.Cells.Item(42, 1).Value = 45796.41666666666
$ws.Cells.Item(42, 2).Value = 5542
$ws.Cells.Item(42, 4).Value = "19.05.202541"
$ws.Cells.Item(43, 1).Value = 45796.42708333334
$ws.Cells.Item(43, 2).Value = 5432
$ws.Cells.Item(43, 4).Value = "19.05.202542"
$ws.Cells.Item(44, 1).Value = 45796.4375
$ws.Cells.Item(44, 2).Value = 5435
$ws.Cells.Item(44, 4).Value = "19.05.202543"
$ws.Cells.Item(45, 1).Value = 45796.44791666666
$ws.Cells.Item(45, 2).Value = 5495
$ws.Cells.Item(45, 4).Value = "19.05.202544"
$ws.Cells.Item(46, 1).Value = 45796.45833333334
$ws.Cells.Item(46, 2).Value = 5431
$ws.Cells.Item(46, 4).Value = "19.05.202545"
$ws.Cells.Item(47, 1).Value = 45796.46875
$ws.Cells.Item(47, 2).Value = 5399
$ws.Cells.Item(47, 4).Value = "19.05.202546"
$ws.Cells.Item(48, 1).Value = 45796.47916666666
$ws.Cells.Item(48, 2).Value = 5350
$ws.Cells.Item(48, 4).Value = "19.05.202547"
$ws.Cells.Item(49, 1).Value = 45796.48958333334
$ws.Cells.Item(49, 2).Value = 5353
$ws.Cells.Item(49, 4).Value = "19.05.202548"
$ws.Cells.Item(50, 1).Value = 45796.5
$ws.Cells.Item(50, 2).Value = 5354
$ws.Cells.Item(50, 4).Value = "19.05.202549"
$ws.Cells.Item(51, 1).Value = 45796.51041666666
$ws.Cells.Item(51, 2).Value = 5361
$ws.Cells.Item(51, 4).Value = "19.05.202550"
$ws.Cells.Item(52, 1).Value = 45796.52083333334
$ws.Cells.Item(52, 2).Value = 5375
$ws.Cells.Item(52, 4).Value = "19.05.202551"
$ws.Cells.Item(53, 1).Value = 45796.53125
$ws.Cells.Item(53, 2).Value = 5422
$ws.Cells.Item(53, 4).Value = "19.05.202552"
$ws.Cells.Item(54, 1).Value = 45796.54166666666
$ws.Cells.Item(54, 2).Value = 5299
$ws.Cells.Item(54, 4).Value = "19.05.202553"
$ws.Cells.Item(55, 1).Value = 45796.55208333334
$ws.Cells.Item(55, 2).Value = 5241
$ws.Cells.Item(55, 4).Value = "19.05.202554"
$ws.Cells.Item(56, 1).Value = 45796.5625
$ws.Cells.Item(56, 2).Value = 5291
$ws.Cells.Item(56, 4).Value = "19.05.202555"
$ws.Cells.Item(57, 1).Value = 45796.57291666666
$ws.Cells.Item(57, 2).Value = 5218
$ws.Cells.Item(57, 4).Value = "19.05.202556"
$ws.Cells.Item(58, 1).Value = 45796.58333333334
$ws.Cells.Item(58, 2).Value = 5253
$ws.Cells.Item(58, 4).Value = "19.05.202557"
$ws.Cells.Item(59, 1).Value = 45796.59375
$ws.Cells.Item(59, 2).Value = 5382
$ws.Cells.Item(59, 4).Value = "19.05.202558"
$ws.Cells.Item(60, 1).Value = 45796.60416666666
$ws.Cells.Item(60, 2).Value = 5312
$ws.Cells.Item(60, 4).Value = "19.05.202559"
$ws.Cells.Item(61, 1).Value = 45796.61458333334
$ws.Cells.Item(61, 2).Value = 5400
$ws.Cells.Item(61, 4).Value = "19.05.202560"
$ws.Cells.Item(62, 1).Value = 45796.625
$ws.Cells.Item(62, 2).Value = 5556
$ws.Cells.Item(62, 4).Value = "19.05.202561"
$ws.Cells.Item(63, 1).Value = 45796.63541666666
$ws.Cells.Item(63, 2).Value = 5693
$ws.Cells.Item(63, 4).Value = "19.05.202562"
$ws.Cells.Item(64, 1).Value = 45796.64583333334
$ws.Cells.Item(64, 2).Value = 5787
$ws.Cells.Item(64, 4).Value = "19.05.202563"
$ws.Cells.Item(65, 1).Value = 45796.65625
$ws.Cells.Item(65, 2).Value = 5915
$ws.Cells.Item(65, 4).Value = "19.05.202564"
$ws.Cells.Item(66, 1).Value = 45796.66666666666
$ws.Cells.Item(66, 2).Value = 5878
$ws.Cells.Item(66, 4).Value = "19.05.202565"
$ws.Cells.Item(67, 1).Value = 45796.67708333334
$ws.Cells.Item(67, 2).Value = 5929
$ws.Cells.Item(67, 4).Value = "19.05.202566"
$ws.Cells.Item(68, 1).Value = 45796.6875
$ws.Cells.Item(68, 2).Value = 5970
$ws.Cells.Item(68, 4).Value = "19.05.202567"
$ws.Cells.Item(69, 1).Value = 45796.69791666666
$ws.Cells.Item(69, 2).Value = 5953
$ws.Cells.Item(69, 4).Value = "19.05.202568"
$ws.Cells.Item(70, 1).Value = 45796.70833333334
$ws.Cells.Item(70, 2).Value = 6166
$ws.Cells.Item(70, 4).Value = "19.05.202569"
$ws.Cells.Item(71, 1).Value = 45796.71875
$ws.Cells.Item(71, 2).Value = 6180
$ws.Cells.Item(71, 4).Value = "19.05.202570"
$ws.Cells.Item(72, 1).Value = 45796.72916666666
$ws.Cells.Item(72, 2).Value = 6244
$ws.Cells.Item(72, 4).Value = "19.05.202571"
$ws.Cells.Item(73, 1).Value = 45796.73958333334
$ws.Cells.Item(73, 2).Value = 6286
$ws.Cells.Item(73, 4).Value = "19.05.202572"
$ws.Cells.Item(74, 1).Value = 45796.75
$ws.Cells.Item(74, 2).Value = 6439
$ws.Cells.Item(74, 4).Value = "19.05.202573"
$ws.Cells.Item(75, 1).Value = 45796.76041666666
$ws.Cells.Item(75, 2).Value = 6585
$ws.Cells.Item(75, 4).Value = "19.05.202574"
$ws.Cells.Item(76, 1).Value = 45796.77083333334
$ws.Cells.Item(76, 2).Value = 6596
$ws.Cells.Item(76, 4).Value = "19.05.202575"
$ws.Cells.Item(77, 1).Value = 45796.78125
$ws.Cells.Item(77, 2).Value = 6744
$ws.Cells.Item(77, 4).Value = "19.05.202576"
$ws.Cells.Item(78, 1).Value = 45796.79166666666
$ws.Cells.Item(78, 2).Value = 6789
$ws.Cells.Item(78, 4).Value = "19.05.202577"
$ws.Cells.Item(79, 1).Value = 45796.80208333334
$ws.Cells.Item(79, 2).Value = 6843
$ws.Cells.Item(79, 4).Value = "19.05.202578"
$ws.Cells.Item(80, 1).Value = 45796.8125
$ws.Cells.Item(80, 2).Value = 6940
$ws.Cells.Item(80, 4).Value = "19.05.202579"
$ws.Cells.Item(81, 1).Value = 45796.82291666666
$ws.Cells.Item(81, 2).Value = 7047
$ws.Cells.Item(81, 4).Value = "19.05.202580"
$ws.Cells.Item(82, 1).Value = 45796.83333333334
$ws.Cells.Item(82, 2).Value = 7076
$ws.Cells.Item(82, 4).Value = "19.05.202581"
$ws.Cells.Item(83, 1).Value = 45796.84375
$ws.Cells.Item(83, 2).Value = 7057
$ws.Cells.Item(83, 4).Value = "19.05.202582"
$ws.Cells.Item(84, 1).Value = 45796.85416666666
$ws.Cells.Item(84, 2).Value = 7055
$ws.Cells.Item(84, 4).Value = "19.05.202583"
$ws.Cells.Item(85, 1).Value = 45796.86458333334
$ws.Cells.Item(85, 2).Value = 6912
$ws.Cells.Item(85, 4).Value = "19.05.202584"
$ws.Cells.Item(86, 1).Value = 45796.875
$ws.Cells.Item(86, 2).Value = 6766
$ws.Cells.Item(86, 4).Value = "19.05.202585"
$ws.Cells.Item(87, 1).Value = 45796.88541666666
$ws.Cells.Item(87, 2).Value = 6555
$ws.Cells.Item(87, 4).Value = "19.05.202586"
$ws.Cells.Item(88, 1).Value = 45796.89583333334
$ws.Cells.Item(88, 2).Value = 6409
$ws.Cells.Item(88, 4).Value = "19.05.202587"
$ws.Cells.Item(89, 1).Value = 45796.90625
$ws.Cells.Item(89, 2).Value = 6265
$ws.Cells.Item(89, 4).Value = "19.05.202588"
$ws.Cells.Item(90, 1).Value = 45796.91666666666
$ws.Cells.Item(90, 2).Value = 6079
$ws.Cells.Item(90, 4).Value = "19.05.202589"
$ws.Cells.Item(91, 1).Value = 45796.92708333334
$ws.Cells.Item(91, 2).Value = 5944
$ws.Cells.Item(91, 4).Value = "19.05.202590"
$ws.Cells.Item(92, 1).Value = 45796.9375
$ws.Cells.Item(92, 2).Value = 5831
$ws.Cells.Item(92, 4).Value = "19.05.202591"
$ws.Cells.Item(93, 1).Value = 45796.94791666666
$ws.Cells.Item(93, 2).Value = 5762
$ws.Cells.Item(93, 4).Value = "19.05.202592"
$ws.Cells.Item(94, 1).Value = 45796.95833333334
$ws.Cells.Item(94, 2).Value = 5676
$ws.Cells.Item(94, 4).Value = "19.05.202593"
$ws.Cells.Item(95, 1).Value = 45796.96875
$ws.Cells.Item(95, 2).Value = 5625
$ws.Cells.Item(95, 4).Value = "19.05.202594"
$ws.Cells.Item(96, 1).Value = 45796.97916666666
$ws.Cells.Item(96, 2).Value = 5481
$ws.Cells.Item(96, 4).Value = "19.05.202595"
$ws.Cells.Item(97, 1).Value = 45796.98958333334
$ws.Cells.Item(97, 2).Value = 5430
$ws.Cells.Item(97, 4).Value = "19.05.202596"
$ws.Cells.Item(98, 1).Value = 45797
$ws.Cells.Item(98, 2).Value = 5376
$ws.Cells.Item(98, 4).Value = "20.05.20251"
$ws.Cells.Item(99, 1).Value = 45797.01041666666
$ws.Cells.Item(99, 2).Value = 5355
$ws.Cells.Item(99, 4).Value = "20.05.20252"
$ws.Cells.Item(100, 1).Value = 45797.02083333334
$ws.Cells.Item(100, 2).Value = 5348
$ws.Cells.Item(100, 4).Value = "20.05.20253"
$ws.Cells.Item(101, 1).Value = 45797.03125
$ws.Cells.Item(101, 2).Value = 5284
$ws.Cells.Item(101, 4).Value = "20.05.20254"
$ws.Cells.Item(102, 1).Value = 45797.04166666666
$ws.Cells.Item(102, 2).Value = 5230
$ws.Cells.Item(102, 4).Value = "20.05.20255"
$ws.Cells.Item(103, 1).Value = 45797.05208333334
$ws.Cells.Item(103, 2).Value = 5221
$ws.Cells.Item(103, 4).Value = "20.05.20256"
$ws.Cells.Item(104, 1).Value = 45797.0625
$ws.Cells.Item(104, 2).Value = 5166
$ws.Cells.Item(104, 4).Value = "20.05.20257"
$ws.Cells.Item(105, 1).Value = 45797.07291666666
$ws.Cells.Item(105, 2).Value = 5197
$ws.Cells.Item(105, 4).Value = "20.05.20258"
$ws.Cells.Item(106, 1).Value = 45797.08333333334
$ws.Cells.Item(106, 2).Value = 5159
$ws.Cells.Item(106, 4).Value = "20.05.20259"
$ws.Cells.Item(107, 1).Value = 45797.09375
$ws.Cells.Item(107, 2).Value = 5131
$ws.Cells.Item(107, 4).Value = "20.05.202510"
$ws.Cells.Item(108, 1).Value = 45797.10416666666
$ws.Cells.Item(108, 2).Value = 5124
$ws.Cells.Item(108, 4).Value = "20.05.202511"
$ws.Cells.Item(109, 1).Value = 45797.11458333334
$ws.Cells.Item(109, 2).Value = 5163
$ws.Cells.Item(109, 4).Value = "20.05.202512"
$ws.Cells.Item(110, 1).Value = 45797.125
$ws.Cells.Item(110, 2).Value = 5247
$ws.Cells.Item(110, 4).Value = "20.05.202513"
$ws.Cells.Item(111, 1).Value = 45797.13541666666
$ws.Cells.Item(111, 2).Value = 5190
$ws.Cells.Item(111, 4).Value = "20.05.202514"
$ws.Cells.Item(112, 1).Value = 45797.14583333334
$ws.Cells.Item(112, 2).Value = 5249
$ws.Cells.Item(112, 4).Value = "20.05.202515"
$ws.Cells.Item(113, 1).Value = 45797.15625
$ws.Cells.Item(113, 2).Value = 5234
$ws.Cells.Item(113, 4).Value = "20.05.202516"
$ws.Cells.Item(114, 1).Value = 45797.16666666666
$ws.Cells.Item(114, 2).Value = 5337
$ws.Cells.Item(114, 4).Value = "20.05.202517"
$ws.Cells.Item(115, 1).Value = 45797.17708333334
$ws.Cells.Item(115, 2).Value = 5309
$ws.Cells.Item(115, 4).Value = "20.05.202518"
$ws.Cells.Item(116, 1).Value = 45797.1875
$ws.Cells.Item(116, 2).Value = 5342
$ws.Cells.Item(116, 4).Value = "20.05.202519"
$ws.Cells.Item(117, 1).Value = 45797.19791666666
$ws.Cells.Item(117, 2).Value = 5347
$ws.Cells.Item(117, 4).Value = "20.05.202520"
$ws.Cells.Item(118, 1).Value = 45797.20833333334
$ws.Cells.Item(118, 2).Value = 5505
$ws.Cells.Item(118, 4).Value = "20.05.202521"
$ws.Cells.Item(119, 1).Value = 45797.21875
$ws.Cells.Item(119, 2).Value = 5582
$ws.Cells.Item(119, 4).Value = "20.05.202522"
$ws.Cells.Item(120, 1).Value = 45797.22916666666
$ws.Cells.Item(120, 2).Value = 5686
$ws.Cells.Item(120, 4).Value = "20.05.202523"
$ws.Cells.Item(121, 1).Value = 45797.23958333334
$ws.Cells.Item(121, 2).Value = 5735
$ws.Cells.Item(121, 4).Value = "20.05.202524"
$ws.Cells.Item(122, 1).Value = 45797.25
$ws.Cells.Item(122, 2).Value = 5920
$ws.Cells.Item(122, 4).Value = "20.05.202525"
$ws.Cells.Item(123, 1).Value = 45797.26041666666
$ws.Cells.Item(123, 2).Value = 5975
$ws.Cells.Item(123, 4).Value = "20.05.202526"
$ws.Cells.Item(124, 1).Value = 45797.27083333334
$ws.Cells.Item(124, 2).Value = 6003
$ws.Cells.Item(124, 4).Value = "20.05.202527"
$ws.Cells.Item(125, 1).Value = 45797.28125
$ws.Cells.Item(125, 2).Value = 5979
$ws.Cells.Item(125, 4).Value = "20.05.202528"
$ws.Cells.Item(126, 1).Value = 45797.29166666666
$ws.Cells.Item(126, 2).Value = 6017
$ws.Cells.Item(126, 4).Value = "20.05.202529"
$ws.Cells.Item(127, 1).Value = 45797.30208333334
$ws.Cells.Item(127, 2).Value = 6000
$ws.Cells.Item(127, 4).Value = "20.05.202530"
$ws.Cells.Item(128, 1).Value = 45797.3125
$ws.Cells.Item(128, 2).Value = 5899
$ws.Cells.Item(128, 4).Value = "20.05.202531"
$ws.Cells.Item(129, 1).Value = 45797.32291666666
$ws.Cells.Item(129, 2).Value = 5807
$ws.Cells.Item(129, 4).Value = "20.05.202532"
$ws.Cells.Item(130, 1).Value = 45797.33333333334
$ws.Cells.Item(130, 2).Value = 5652
$ws.Cells.Item(130, 4).Value = "20.05.202533"
$ws.Cells.Item(131, 1).Value = 45797.34375
$ws.Cells.Item(131, 2).Value = 5556
$ws.Cells.Item(131, 4).Value = "20.05.202534"
$ws.Cells.Item(132, 1).Value = 45797.35416666666
$ws.Cells.Item(132, 2).Value = 5471
$ws.Cells.Item(132, 4).Value = "20.05.202535"
$ws.Cells.Item(133, 1).Value = 45797.36458333334
$ws.Cells.Item(133, 2).Value = 5371
$ws.Cells.Item(133, 4).Value = "20.05.202536"
$ws.Cells.Item(134, 1).Value = 45797.375
$ws.Cells.Item(134, 2).Value = 5346
$ws.Cells.Item(134, 4).Value = "20.05.202537"
$ws.Cells.Item(135, 1).Value = 45797.38541666666
$ws.Cells.Item(135, 2).Value = 5320
$ws.Cells.Item(135, 4).Value = "20.05.202538"
$ws.Cells.Item(136, 1).Value = 45797.39583333334
$ws.Cells.Item(136, 2).Value = 5220
$ws.Cells.Item(136, 4).Value = "20.05.202539"
$ws.Cells.Item(137, 1).Value = 45797.40625
$ws.Cells.Item(137, 2).Value = 0
$ws.Cells.Item(137, 4).Value = "20.05.202540"
$ws.Cells.Item(138, 1).Value = 45797.41666666666
$ws.Cells.Item(138, 2).Value = 0
$ws.Cells.Item(138, 4).Value = "20.05.202541"
$ws.Cells.Item(139, 1).Value = 45797.42708333334
$ws.Cells.Item(139, 2).Value = 0
$ws.Cells.Item(139, 4).Value = "20.05.202542"
$ws.Cells.Item(140, 1).Value = 45797.4375
$ws.Cells.Item(140, 2).Value = 0
$ws.Cells.Item(140, 4).Value = "20.05.202543"
$ws.Cells.Item(141, 1).Value = 45797.44791666666
$ws.Cells.Item(141, 2).Value = 0
$ws.Cells.Item(141, 4).Value = "20.05.202544"
$ws.Cells.Item(142, 1).Value = 45797.45833333334
$ws.Cells.Item(142, 2).Value = 0
$ws.Cells.Item(142, 4).Value = "20.05.202545"
$ws.Cells.Item(143, 1).Value = 45797.46875
$ws.Cells.Item(143, 2).Value = 0
$ws.Cells.Item(143, 4).Value = "20.05.202546"
$ws.Cells.Item(144, 1).Value = 45797.47916666666
$ws.Cells.Item(144, 2).Value = 0
$ws.Cells.Item(144, 4).Value = "20.05.202547"
$ws.Cells.Item(145, 1).Value = 45797.48958333334
$ws.Cells.Item(145, 2).Value = 0
$ws.Cells.Item(145, 4).Value = "20.05.202548"
$ws.Cells.Item(146, 1).Value = 45797.5
$ws.Cells.Item(146, 2).Value = 0
$ws.Cells.Item(146, 4).Value = "20.05.202549"
$ws.Cells.Item(147, 1).Value = 45797.51041666666
$ws.Cells.Item(147, 2).Value = 0
$ws.Cells.Item(147, 4).Value = "20.05.202550"
$ws.Cells.Item(148, 1).Value = 45797.52083333334
$ws.Cells.Item(148, 2).Value = 0
$ws.Cells.Item(148, 4).Value = "20.05.202551"
$ws.Cells.Item(149, 1).Value = 45797.53125
$ws.Cells.Item(149, 2).Value = 0
$ws.Cells.Item(149, 4).Value = "20.05.202552"
$ws.Cells.Item(150, 1).Value = 45797.54166666666
$ws.Cells.Item(150, 2).Value = 0
$ws.Cells.Item(150, 4).Value = "20.05.202553"
$ws.Cells.Item(151, 1).Value = 45797.55208333334
$ws.Cells.Item(151, 2).Value = 0
$ws.Cells.Item(151, 4).Value = "20.05.202554"
$ws.Cells.Item(152, 1).Value = 45797.5625
$ws.Cells.Item(152, 2).Value = 0
$ws.Cells.Item(152, 4).Value = "20.05.202555"
$ws.Cells.Item(153, 1).Value = 45797.57291666666
$ws.Cells.Item(153, 2).Value = 0
$ws.Cells.Item(153, 4).Value = "20.05.202556"
$ws.Cells.Item(154, 1).Value = 45797.58333333334
$ws.Cells.Item(154, 2).Value = 0
$ws.Cells.Item(154, 4).Value = "20.05.202557"
$ws.Cells.Item(155, 1).Value = 45797.59375
$ws.Cells.Item(155, 2).Value = 0
$ws.Cells.Item(155, 4).Value = "20.05.202558"
$ws.Cells.Item(156, 1).Value = 45797.60416666666
$ws.Cells.Item(156, 2).Value = 0
$ws.Cells.Item(156, 4).Value = "20.05.202559"
$ws.Cells.Item(157, 1).Value = 45797.61458333334
$ws.Cells.Item(157, 2).Value = 0
$ws.Cells.Item(157, 4).Value = "20.05.202560"
$ws.Cells.Item(158, 1).Value = 45797.625
$ws.Cells.Item(158, 2).Value = 0
$ws.Cells.Item(158, 4).Value = "20.05.202561"
$ws.Cells.Item(159, 1).Value = 45797.63541666666
$ws.Cells.Item(159, 2).Value = 0
$ws.Cells.Item(159, 4).Value = "20.05.202562"
$ws.Cells.Item(160, 1).Value = 45797.64583333334
$ws.Cells.Item(160, 2).Value = 0
$ws.Cells.Item(160, 4).Value = "20.05.202563"
$ws.Cells.Item(161, 1).Value = 45797.65625
$ws.Cells.Item(161, 2).Value = 0
$ws.Cells.Item(161, 4).Value = "20.05.202564"
$ws.Cells.Item(162, 1).Value = 45797.66666666666
$ws.Cells.Item(162, 2).Value = 0
$ws.Cells.Item(162, 4).Value = "20.05.202565"
$ws.Cells.Item(163, 1).Value = 45797.67708333334
$ws.Cells.Item(163, 2).Value = 0
$ws.Cells.Item(163, 4).Value = "20.05.202566"
$ws.Cells.Item(164, 1).Value = 45797.6875
$ws.Cells.Item(164, 2).Value = 0
$ws.Cells.Item(164, 4).Value = "20.05.202567"
$ws.Cells.Item(165, 1).Value = 45797.69791666666
$ws.Cells.Item(165, 2).Value = 0
$ws.Cells.Item(165, 4).Value = "20.05.202568"
$ws.Cells.Item(166, 1).Value = 45797.70833333334
$ws.Cells.Item(166, 2).Value = 0
$ws.Cells.Item(166, 4).Value = "20.05.202569"
$ws.Cells.Item(167, 1).Value = 45797.71875
$ws.Cells.Item(167, 2).Value = 0
$ws.Cells.Item(167, 4).Value = "20.05.202570"
$ws.Cells.Item(168, 1).Value = 45797.72916666666
$ws.Cells.Item(168, 2).Value = 0
$ws.Cells.Item(168, 4).Value = "20.05.202571"
$ws.Cells.Item(169, 1).Value = 45797.73958333334
$ws.Cells.Item(169, 2).Value = 0
$ws.Cells.Item(169, 4).Value = "20.05.202572"
$ws.Cells.Item(170, 1).Value = 45797.75
$ws.Cells.Item(170, 2).Value = 0
$ws.Cells.Item(170, 4).Value = "20.05.202573"
$ws.Cells.Item(171, 1).Value = 45797.76041666666
$ws.Cells.Item(171, 2).Value = 0
$ws.Cells.Item(171, 4).Value = "20.05.202574"
$ws.Cells.Item(172, 1).Value = 45797.77083333334
$ws.Cells.Item(172, 2).Value = 0
$ws.Cells.Item(172, 4).Value = "20.05.202575"
$ws.Cells.Item(173, 1).Value = 45797.78125
$ws.Cells.Item(173, 2).Value = 0
$ws.Cells.Item(173, 4).Value = "20.05.202576"
$ws.Cells.Item(174, 1).Value = 45797.79166666666
$ws.Cells.Item(174, 2).Value = 0
$ws.Cells.Item(174, 4).Value = "20.05.202577"
$ws.Cells.Item(175, 1).Value = 45797.80208333334
$ws.Cells.Item(175, 2).Value = 0
$ws.Cells.Item(175, 4).Value = "20.05.202578"
$ws.Cells.Item(176, 1).Value = 45797.8125
$ws.Cells.Item(176, 2).Value = 0
$ws.Cells.Item(176, 4).Value = "20.05.202579"
$ws.Cells.Item(177, 1).Value = 45797.82291666666
$ws.Cells.Item(177, 2).Value = 0
$ws.Cells.Item(177, 4).Value = "20.05.202580"
$ws.Cells.Item(178, 1).Value = 45797.83333333334
$ws.Cells.Item(178, 2).Value = 0
$ws.Cells.Item(178, 4).Value = "20.05.202581"
$ws.Cells.Item(179, 1).Value = 45797.84375
$ws.Cells.Item(179, 2).Value = 0
$ws.Cells.Item(179, 4).Value = "20.05.202582"
$ws.Cells.Item(180, 1).Value = 45797.85416666666
$ws.Cells.Item(180, 2).Value = 0
$ws.Cells.Item(180, 4).Value = "20.05.202583"
$ws.Cells.Item(181, 1).Value = 45797.86458333334
$ws.Cells.Item(181, 2).Value = 0
$ws.Cells.Item(181, 4).Value = "20.05.202584"
$ws.Cells.Item(182, 1).Value = 45797.875
$ws.Cells.Item(182, 2).Value = 0
$ws.Cells.Item(182, 4).Value = "20.05.202585"
$ws.Cells.Item(183, 1).Value = 45797.88541666666
$ws.Cells.Item(183, 2).Value = 0
$ws.Cells.Item(183, 4).Value = "20.05.202586"
$ws.Cells.Item(184, 1).Value = 45797.89583333334
$ws.Cells.Item(184, 2).Value = 0
$ws.Cells.Item(184, 4).Value = "20.05.202587"
$ws.Cells.Item(185, 1).Value = 45797.90625
$ws.Cells.Item(185, 2).Value = 0
$ws.Cells.Item(185, 4).Value = "20.05.202588"
$ws.Cells.Item(186, 1).Value = 45797.91666666666
$ws.Cells.Item(186, 2).Value = 0
$ws.Cells.Item(186, 4).Value = "20.05.202589"
$ws.Cells.Item(187, 1).Value = 45797.92708333334
$ws.Cells.Item(187, 2).Value = 0
$ws.Cells.Item(187, 4).Value = "20.05.202590"
$ws.Cells.Item(188, 1).Value = 45797.9375
$ws.Cells.Item(188, 2).Value = 0
$ws.Cells.Item(188, 4).Value = "20.05.202591"
$ws.Cells.Item(189, 1).Value = 45797.94791666666
$ws.Cells.Item(189, 2).Value = 0
$ws.Cells.Item(189, 4).Value = "20.05.202592"
$ws.Cells.Item(190, 1).Value = 45797.95833333334
$ws.Cells.Item(190, 2).Value = 0
$ws.Cells.Item(190, 4).Value = "20.05.202593"
$ws.Cells.Item(191, 1).Value = 45797.96875
$ws.Cells.Item(191, 2).Value = 0
$ws.Cells.Item(191, 4).Value = "20.05.202594"
$ws.Cells.Item(192, 1).Value = 45797.97916666666
$ws.Cells.Item(192, 2).Value = 0
$ws.Cells.Item(192, 4).Value = "20.05.202595"
$ws.Cells.Item(193, 1).Value = 45797.98958333334
$ws.Cells.Item(193, 2).Value = 0
$ws.Cells.Item(193, 4).Value = "20.05.202596"
